# Insert a new weekly price record at row 57 for "Feria Lagunitas de Puerto
# Montt" - Mandarina, shifting the existing rows 57:93 down to 58:94.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:57").Insert()

$ws.Range("A57").Value = 4
$ws.Range("B57").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C57").Value = "Los Lagos"
$ws.Range("D57").Value = 44438
$ws.Range("E57").Value = 10
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100102
$ws.Range("H57").Value = "Cítricos"
$ws.Range("I57").Value = 100102004
$ws.Range("J57").Value = "Mandarina"
$ws.Range("K57").Value = "Clementina"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 300
$ws.Range("N57").Value = 7000
$ws.Range("O57").Value = 7000
$ws.Range("P57").Value = 7000
$ws.Range("Q57").Value = "$/bandeja 10 kilos"
$ws.Range("R57").Value = "Provincia de Limarí"
$ws.Range("S57").Value = 700
$ws.Range("T57").Value = 10
